$p = $ppt.ActivePresentation

# Slide 2, content placeholder shape holds the "Mail: $developers.Name" line.
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Second paragraph (lvl=1): "Mail: $developers.Name"
$para = $tr.Paragraphs(2)

# Rename the merge-field from developers.Name to developers.Mail
# ("Mail: $" is 7 chars, "developers.Name" is the next 15 chars).
$field = $para.Characters(8, 15)
$field.Text = "developers.Mail"

# Turn "$developers.Mail" (16 chars, starting right after "Mail: ") into a
# mailto hyperlink, matching the new <a:hlinkClick r:id="rId2"/> runs.
$link = $para.Characters(7, 16)
$actionSetting = $link.ActionSettings(1)
$actionSetting.Hyperlink.Address = "mailto:"
